$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column updates ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 475
$ws1.Range("F4").Value = 8002
$ws1.Range("F8").Value = 31
$ws1.Range("F13").Value = 456
$ws1.Range("F16").Value = 30
$ws1.Range("F17").Value = 5895
$ws1.Range("F18").Value = 186
$ws1.Range("F19").Value = 271
$ws1.Range("F20").Value = 1888
$ws1.Range("F21").Value = 11
$ws1.Range("F22").Value = 23
$ws1.Range("F24").Value = 404

# Sheet "全部类型" (sheet4) - F column updates ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 475
$ws4.Range("F4").Value = 8002
$ws4.Range("F8").Value = 31
$ws4.Range("F13").Value = 456
$ws4.Range("F16").Value = 30
$ws4.Range("F18").Value = 5895
$ws4.Range("F20").Value = 186
$ws4.Range("F21").Value = 271
$ws4.Range("F22").Value = 1888
$ws4.Range("F23").Value = 11
$ws4.Range("F24").Value = 23
$ws4.Range("F26").Value = 404

$wb.Save()
